$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: write a brand-new "data" row (columns A-G) that mirrors the look of
# an existing reference row ($styleRow), using date-like text in column A
# (forced to Text via NumberFormat "@" so Excel doesn't auto-convert it to a
# serial date), a real datetime serial in column B (copying the reference
# row's datetime format), and plain numbers in D:G.
# ---------------------------------------------------------------------------
function Add-DataRow {
    param(
        $ws,
        [int]$row,
        [int]$styleRow,
        [string]$dateText,
        [double]$timeVal,
        [string]$name,
        $d, $e, $f, $g,
        [bool]$applyStyle = $true
    )

    $ws.Range("A$row").NumberFormat = "@"
    $ws.Range("A$row").Value = $dateText
    $ws.Range("B$row").Value = $timeVal
    $ws.Range("C$row").Value = $name
    $ws.Range("D$row").Value = $d
    $ws.Range("E$row").Value = $e
    $ws.Range("F$row").Value = $f
    $ws.Range("G$row").Value = $g

    # Column B always keeps the datetime display format of the reference row.
    $ws.Range("B$styleRow").Copy()
    $ws.Range("B$row").PasteSpecial(-4122)

    if ($applyStyle) {
        $ws.Range("A$row").Style = $ws.Range("A$styleRow").Style
        $ws.Range("C$row").Style = $ws.Range("C$styleRow").Style
        $ws.Range("D$row").Style = $ws.Range("D$styleRow").Style
        $ws.Range("E$row").Style = $ws.Range("E$styleRow").Style
        $ws.Range("F$row").Style = $ws.Range("F$styleRow").Style
        $ws.Range("G$row").Style = $ws.Range("G$styleRow").Style
    }
}

# ---------------------------------------------------------------------------
# Sheet "AMSIN": append sprint-run history rows 80-84
# ---------------------------------------------------------------------------
$wsAmsin = $wb.Worksheets.Item("AMSIN")

Add-DataRow $wsAmsin 80 79 "2024-03-28" 45379.52823976852 "190fstocrr" 41 41 0 1.77
Add-DataRow $wsAmsin 81 79 "2024-03-29" 45380.47474833333 "190ocrscnd" 41 41 0 1.44
Add-DataRow $wsAmsin 82 79 "2024-04-01" 45383.36676361111 "190fnlocrr" 41 41 0 1.44
Add-DataRow $wsAmsin 83 79 "2024-05-02" 45414.47038731482 "191fstocr"  41 41 0 1.89
Add-DataRow $wsAmsin 84 79 "2024-05-03" 45415.34386107639 "191lstocr"  41 41 0 1.51

# ---------------------------------------------------------------------------
# Sheet "AMS": re-run/fix row 64 (bug in the 191 series uncovered a stale
# timestamp + missing formatting on the previous row), then append rows
# 65-67. Row 67 (the newest "191betaocr" entry) keeps the legacy/unstyled
# look, matching how the tracker had originally wired in row 64 before the
# fix below.
# ---------------------------------------------------------------------------
$wsAms = $wb.Worksheets.Item("AMS")

# Fix the stale timestamp on row 64 and normalize its formatting to match
# the rest of the table (copy the look of row 63, which was never touched).
$wsAms.Range("B64").Value = 45364.63885644676
$wsAms.Range("A64").Style = $wsAms.Range("A63").Style
$wsAms.Range("C64").Style = $wsAms.Range("C63").Style
$wsAms.Range("D64").Style = $wsAms.Range("D63").Style
$wsAms.Range("E64").Style = $wsAms.Range("E63").Style
$wsAms.Range("F64").Style = $wsAms.Range("F63").Style
$wsAms.Range("G64").Style = $wsAms.Range("G63").Style

Add-DataRow $wsAms 65 64 "2024-04-01" 45383.52811969908 "190betaocr" 41 41 0 1.28
Add-DataRow $wsAms 66 64 "2024-04-01" 45383.84961586806 "190liveocr" 41 41 0 1.28
Add-DataRow $wsAms 67 64 "2024-05-03" 45415.60354174538 "191betaocr" 41 41 0 1.19 $false
